$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("C").Delete()
$ws.Range("C1").Select()
